$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-520)
# from serial date 45205 to 45206 (2023-10-06 -> 2023-10-07)
$ws.Range("C2:C520").Value = 45206
